# Extend the "Path to Graduation" schedule with additional semester blocks
# (Fall 2024/Spring 2024/Summer 2024 course rows, plus brand-new
# Fall 2025/Spring 2025/Summer 2025 and Fall 2026/Spring 2026/Summer 2026
# blocks), and rebalance the existing Fall 2023/Spring 2023/Summer 2023
# course rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fall 2022 / Spring 2022 / Summer 2022 block: drop the last two course rows ---
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()

# --- Fall 2023 / Spring 2023 / Summer 2023 block: rebalance course rows 13-14 ---
$ws.Range("C13").Value = "CYBR 3106"
$ws.Range("D13").Value = 3
$ws.Range("C14").Value = "CYBR 3108"
$ws.Range("D14").Value = 3
$ws.Range("A15").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("A16").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("A17").ClearContents()
$ws.Range("B17").ClearContents()

# --- Fall 2024 / Spring 2024 / Summer 2024 block: fill in course rows 22-23 ---
$ws.Range("A22").Value = "CYBR 3115"
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = "CPSC 4111"
$ws.Range("D22").Value = 3
$ws.Range("A23").Value = "CYBR 3119"
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = "CPSC 4115"
$ws.Range("D23").Value = 3

# --- New Fall 2025 / Spring 2025 / Summer 2025 block ---
$ws.Range("A30").Value = "Fall 2025"
$ws.Range("B30").Value = "Credits"
$ws.Range("C30").Value = "Spring 2025"
$ws.Range("D30").Value = "Credits"
$ws.Range("E30").Value = "Summer 2025"
$ws.Range("F30").Value = "Credits"

$ws.Range("A31").Value = "CPSC 6180"
$ws.Range("B31").Value = 3
$ws.Range("C31").Value = "CPSC 6985"
$ws.Range("D31").Value = 4

$ws.Range("A32").Value = "CPSC 6185"
$ws.Range("B32").Value = 3

$ws.Range("A38").Value = "Total"
$ws.Range("B38").Formula = "=SUM(B31:B37)"
$ws.Range("C38").Value = "Total"
$ws.Range("D38").Formula = "=SUM(D31:D37)"
$ws.Range("E38").Value = "Total"
$ws.Range("F38").Formula = "=SUM(F31:F37)"

# --- New Fall 2026 / Spring 2026 / Summer 2026 block ---
$ws.Range("A39").Value = "Fall 2026"
$ws.Range("B39").Value = "Credits"
$ws.Range("C39").Value = "Spring 2026"
$ws.Range("D39").Value = "Credits"
$ws.Range("E39").Value = "Summer 2026"
$ws.Range("F39").Value = "Credits"

$ws.Range("A47").Value = "Total"
$ws.Range("B47").Formula = "=SUM(B40:B46)"
$ws.Range("C47").Value = "Total"
$ws.Range("D47").Formula = "=SUM(D40:D46)"
$ws.Range("E47").Value = "Total"
$ws.Range("F47").Formula = "=SUM(F40:F46)"
